$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "1.001", "0.000007124") keep their exact original formatting
# instead of being auto-converted to floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.507.67"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.730.42"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "244.85"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.4917"
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("D8").Value = "0.2636"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("D9").Value = "0.06186"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "1.737.43"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("D11").Value = "0.07025"
$ws.Range("E11").Value = "  -1.73%  "
$ws.Range("D12").Value = "15.58"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "4.551"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "0.6009"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "77.38"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "26.528.58"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "0.000007124"
$ws.Range("E19").Value = "  +2.85%  "
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").Value = "1.965.12"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").Value = "4.498"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "8.585"
$ws.Range("E23").Value = "  -4.05%  "
$ws.Range("D24").Value = "5.190"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "138.57"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").Value = "15.24"
$ws.Range("D27").Value = "1.438"
$ws.Range("E27").Value = "  +2.48%  "
$ws.Range("D28").Value = "106.76"
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("D29").Value = "1.723"
$ws.Range("E29").Value = "  -3.96%  "
$ws.Range("D30").Value = "3.966"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").Value = "0.04528"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "1.001"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6260"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.9074"
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.018"
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.402"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.01491"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "100.71"
$ws.Range("E43").Value = "  -3.81%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.464"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3876"
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "6.693"
$ws.Range("E46").Value = "  -3.16%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1159"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05366"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "30.35"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "7.681"
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.237"
$ws.Range("E51").Value = "  -2.43%  "
